$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "27.729.71"
Set-TextValue "E2" "  -1.72%  "
Set-TextValue "D3" "1.761.60"
Set-TextValue "E3" "  -1.76%  "
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "327.10"
Set-TextValue "E5" "  -1.95%  "
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  -0.11%  "
Set-TextValue "D7" "0.4411"
Set-TextValue "E7" "  -1.76%  "
Set-TextValue "D8" "0.3747"
Set-TextValue "E8" "  +1.26%  "
Set-TextValue "D9" "45.46"
Set-TextValue "E9" "  +1.07%  "
Set-TextValue "D10" "0.07785"
Set-TextValue "E10" "  +3.31%  "
Set-TextValue "D11" "1.128"
Set-TextValue "E11" "  -0.94%  "
Set-TextValue "E12" "  -0.14%  "
Set-TextValue "D13" "21.77"
Set-TextValue "E13" "  -2.26%  "
Set-TextValue "D14" "6.198"
Set-TextValue "E14" "  -1.14%  "
Set-TextValue "D15" "7.402"
Set-TextValue "E15" "  -0.68%  "
Set-TextValue "D16" "1.757.92"
Set-TextValue "E16" "  -2.00%  "
Set-TextValue "B17" "Litecoin"
Set-TextValue "C17" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D17" "91.12"
Set-TextValue "E17" "  +12.80%  "
Set-TextValue "B18" "ShibaInu"
Set-TextValue "C18" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D18" "0.00001084"
Set-TextValue "E18" "  -0.10%  "
Set-TextValue "D19" "0.06243"
Set-TextValue "E19" "  -7.29%  "
Set-TextValue "E20" "  -0.07%  "
Set-TextValue "D21" "17.36"
Set-TextValue "E21" "  -0.20%  "
Set-TextValue "D22" "6.188"
Set-TextValue "E22" "  -2.53%  "
Set-TextValue "D23" "0.5323"
Set-TextValue "E23" "  -2.92%  "
Set-TextValue "D24" "27.755.14"
Set-TextValue "E24" "  -1.55%  "
Set-TextValue "D25" "11.66"
Set-TextValue "E25" "  -0.88%  "
Set-TextValue "D26" "2.313"
Set-TextValue "E26" "  -4.22%  "
Set-TextValue "D27" "20.85"
Set-TextValue "E27" "  +1.71%  "
Set-TextValue "D28" "153.72"
Set-TextValue "E28" "  +1.37%  "
Set-TextValue "D29" "2.358"
Set-TextValue "E29" "  +0.18%  "
Set-TextValue "D30" "1.956.63"
Set-TextValue "E30" "  -2.02%  "
Set-TextValue "D31" "129.11"
Set-TextValue "E31" "  -2.34%  "
Set-TextValue "D32" "1.218"
Set-TextValue "E32" "  -1.35%  "
Set-TextValue "D33" "5.792"
Set-TextValue "E33" "  +0.28%  "
Set-TextValue "D34" "0.09300"
Set-TextValue "E34" "  -0.88%  "
Set-TextValue "D35" "3.644"
Set-TextValue "E35" "  -9.81%  "
Set-TextValue "D36" "12.72"
Set-TextValue "E36" "  +5.77%  "
Set-TextValue "E37" "  -8.60%  "
Set-TextValue "D38" "0.02327"
Set-TextValue "E38" "  -0.31%  "
Set-TextValue "D39" "0.06163"
Set-TextValue "E39" "  -2.10%  "
Set-TextValue "D40" "0.6506"
Set-TextValue "E40" "  -0.46%  "
Set-TextValue "D41" "5.104"
Set-TextValue "E41" "  -1.31%  "
Set-TextValue "D42" "1.201"
Set-TextValue "E42" "  +0.12%  "
Set-TextValue "D43" "8.024"
Set-TextValue "E43" "  -3.51%  "
Set-TextValue "E44" "  -3.98%  "
Set-TextValue "E45" "  -0.05%  "
Set-TextValue "D46" "13.75"
Set-TextValue "E46" "  -2.46%  "
Set-TextValue "D47" "0.6033"
Set-TextValue "E47" "  -0.35%  "
Set-TextValue "D48" "3.763"
Set-TextValue "E48" "  -1.41%  "
Set-TextValue "D49" "126.32"
Set-TextValue "E49" "  -2.46%  "
Set-TextValue "D50" "2.002"
Set-TextValue "E50" "  -1.04%  "
Set-TextValue "D51" "1.146"
Set-TextValue "E51" "  -0.73%  "
